# Fix word order in the "Härkä tähtikuvio" (Taurus constellation) heading.
# Old: "havainnointijaksot vuonna Härkä tähtikuvio 2022: tammikuuta 16-25"
# New: "Härkä tähtikuvio havainnointijaksot vuonna 2022: tammikuuta 16-25"
# This phrase occurs 4 times in the document, all identical, so a global
# Find/Replace (Wrap=2 i.e. wdFindContinue, Replace:=2 i.e. wdReplaceAll)
# over the whole document content covers every occurrence.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "havainnointijaksot vuonna Härkä tähtikuvio 2022: tammikuuta 16-25", `
    $true, `
    $false, `
    $false, `
    $false, `
    $false, `
    $true, `
    1, `
    $false, `
    "Härkä tähtikuvio havainnointijaksot vuonna 2022: tammikuuta 16-25", `
    2
)

$d.Save()
